$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("A2").Value = "https://stackoverflow.com/questions/50815522/selenium-python-automation?r=SearchResults"
$ws.Range("A3").Value = "https://stackoverflow.com/questions/tagged/python"
$ws.Range("A4").Value = "https://stackoverflow.com/questions/tagged/selenium"
$ws.Range("A5").Value = "https://stackoverflow.com/questions/tagged/web-crawler"
$ws.Range("A6").Value = "https://stackoverflow.com/users/9834021/vidhya"
$ws.Range("A7").Value = "https://stackoverflow.com/users/9834021/vidhya"
$ws.Range("A8").Value = "https://stackoverflow.com/questions/43637687/python-automation?r=SearchResults"
$ws.Range("A9").Value = "https://stackoverflow.com/questions/tagged/python"
$ws.Range("A10").Value = "https://stackoverflow.com/questions/tagged/selenium"
$ws.Range("A11").Value = "https://stackoverflow.com/questions/tagged/automation"
